$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$newQuery = "MATCH (p:participant)-->(s:study)`r`nOPTIONAL MATCH (samp:sample)-->(p)`r`nOPTIONAL MATCH (p)<--(diag:diagnosis)`r`nOPTIONAL MATCH (samp)<--(f:file)`r`nOPTIONAL MATCH (f)<--(g:genomic_info)`r`nWITH s, p, samp, f, g, diag`r`nWHERE g.platform in ['NovaSeqS4']`r`nwith p`r`nOPTIONAL MATCH (p)-->(s:study)`r`nOPTIONAL MATCH (samp:sample)-->(p)`r`nWITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp`r`nRETURN`r`ncoalesce(p.participant_id,'') as ``Participant ID``,`r`ncoalesce(s.study_name, '') as ``Study Name``,`r`ncoalesce(s.phs_accession,'') as ``Accession``,`r`ncoalesce(p.gender,'') as ``Gender``,`r`ncoalesce(apoc.text.join(samp, ','), '') as ``Samples```r`nORDER BY p.participant_id LIMIT 100"

$ws.Range("B2").Value = $newQuery
$ws.Rows.Item(2).RowHeight = 279

$ws.Range("B4").Select()
